# Adding contributors to the presentation
# Adds a "Submitted by:" textbox with the team members' names to the
# title slide (slide 1), placed after the existing picture shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Position/size below (EMU -> points, 1 pt = 1/72 in = 12700 EMU):
#   off  x=7600426  y=4630723   EMU
#   ext  cx=3900881 cy=1200329  EMU
$left   = 598.4587401574803
$top    = 364.6238582677165
$width  = 307.1559906519685
$height = 94.51409538818898

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$tf = $tb.TextFrame
$tf.WordWrap = -1

$tr = $tf.TextRange
$tr.Text = "Submitted by:"
$tr.LanguageID = "en-IN"
[void]$tr.InsertAfter("`nAryan Kumar (M23CSA510)")
[void]$tr.InsertAfter("`nHarsh Parashar (M22AIE210)")
[void]$tr.InsertAfter("`nPrateek Singhal (M22AIE215)")

# Resize the shape height to fit the text (matches <a:spAutoFit/>).
$tf.AutoSize = 1

# No shape fill (<a:noFill/>).
$tb.Fill.Visible = $false

$tb.Name = "TextBox 3"
